$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Global" sheet lists three device rows used for cross-browser/device
# testing. The commit swaps the unstable iOS device (iPhone 6s) for a more
# stable one (iPhone 8, with a new device_id) and moves that row to the
# bottom of the device list, while the two Android device rows move up.

# Row 4 (was row 5): Android device - SM-G950F, no URL, com.Advantage.aShopping
$ws.Cells.Item(4,3).Value = ""
$ws.Cells.Item(4,4).Value = "ce031713bc66a70d05"
$ws.Cells.Item(4,5).Value = "com.Advantage.aShopping"
$ws.Cells.Item(4,6).Value = $false
$ws.Cells.Item(4,7).Value = "ANDROID"
$ws.Cells.Item(4,8).Value = "SM-G950F"

# Row 5 (was row 6): Android device - SM-G950F, with URL, MC.Browser
$ws.Cells.Item(5,3).Value = "https://advantageonlineshopping.com/"
$ws.Cells.Item(5,4).Value = "ce031713bc66a70d05"
$ws.Cells.Item(5,5).Value = "MC.Browser"
$ws.Cells.Item(5,6).Value = $false
$ws.Cells.Item(5,7).Value = "ANDROID"
$ws.Cells.Item(5,8).Value = "SM-G950F"

# Row 6 (was row 4): iOS device - now iPhone 8 (formerly iPhone 6s) with a
# new, more stable device_id, com.hpe.iShopping
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = "ed2ff5276810f2265b87cb2d58acc7b9246aa5c4"
$ws.Cells.Item(6,5).Value = "com.hpe.iShopping"
$ws.Cells.Item(6,6).Value = $true
$ws.Cells.Item(6,7).Value = "IOS"
$ws.Cells.Item(6,8).Value = "iPhone 8"

# The device_id column best-fit width shifts slightly because the new
# device identifier string renders at a different width.
$ws.Columns(4).ColumnWidth = 37.45
